$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (columns D, K, L, M, N, O, P, Q, R, S, T)
# derived from a rotation of the existing data rows 2-6.
$data = @{
    2 = @{ D=44210; K="Rainier"; L="Segunda"; M=250; N=21000; O=22000; P=21500; Q="`$/caja 18 kilos";     R="Región de O'Higgins";  S=1194; T=18 }
    3 = @{ D=44161; K="Bing";    L="Primera"; M=160; N=39000; O=40000; P=39500; Q="`$/caja 20 kilos";     R="Provincia de Curicó";  S=1975; T=20 }
    4 = @{ D=44208; K="Lapins";  L="Segunda"; M=200; N=10500; O=11000; P=10750; Q="`$/bandeja 12 kilos";  R="Provincia de Curicó";  S=896;  T=12 }
    5 = @{ D=44229; K="Santina"; L="Primera"; M=250; N=6500;  O=7000;  P=6750;  Q="`$/bandeja 5 kilos";   R="Provincia de Curicó";  S=1350; T=5  }
    6 = @{ D=44175; K="Rainier"; L="Segunda"; M=270; N=25000; O=26000; P=25500; Q="`$/caja 18 kilos";     R="Región de O'Higgins";  S=1417; T=18 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value  = $vals.D   # D - Fecha
    $ws.Cells.Item($row, 11).Value = $vals.K   # K - Variedad
    $ws.Cells.Item($row, 12).Value = $vals.L   # L - Calidad
    $ws.Cells.Item($row, 13).Value = $vals.M   # M - Volumen
    $ws.Cells.Item($row, 14).Value = $vals.N   # N - Precio minimo
    $ws.Cells.Item($row, 15).Value = $vals.O   # O - Precio maximo
    $ws.Cells.Item($row, 16).Value = $vals.P   # P - Precio promedio ponderado
    $ws.Cells.Item($row, 17).Value = $vals.Q   # Q - Unidad de comercializacion
    $ws.Cells.Item($row, 18).Value = $vals.R   # R - Origen
    $ws.Cells.Item($row, 19).Value = $vals.S   # S - Precio $/Kg
    $ws.Cells.Item($row, 20).Value = $vals.T   # T - Kg / unidad
}
